# Hortaliza, Vega Modelo de Temuco - Repollo
# Insert 2 new weekly rows (665-666) above the former row 665, shifting
# every subsequent row down by two (old 665..711 -> new 667..713).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 665-666; existing rows 665.. shift down to 667..
$ws.Range("A665:R666").Insert()

# Constant columns shared by every data row in this sheet.
$mercadoId   = 10
$mercado     = "Vega Modelo de Temuco"
$region      = "La Araucanía"
$codreg      = 9
$categoriaId = 100112006
$categoria   = "Repollo"
$calidad     = "Primera"
$unidad      = "`$/unidad"
$kgUnidades  = 1
$clasif      = "Hortaliza"

# New row 665: Crespo record
$r = 665
$ws.Cells.Item($r, 1).Value  = $mercadoId
$ws.Cells.Item($r, 2).Value  = $mercado
$ws.Cells.Item($r, 3).Value  = $region
$ws.Cells.Item($r, 4).Value  = 44746
$ws.Cells.Item($r, 5).Value  = $codreg
$ws.Cells.Item($r, 6).Value  = $categoriaId
$ws.Cells.Item($r, 7).Value  = $categoria
$ws.Cells.Item($r, 8).Value  = "Crespo record"
$ws.Cells.Item($r, 9).Value  = $calidad
$ws.Cells.Item($r, 10).Value = 1700
$ws.Cells.Item($r, 11).Value = 1300
$ws.Cells.Item($r, 12).Value = 1500
$ws.Cells.Item($r, 13).Value = 1353
$ws.Cells.Item($r, 14).Value = $unidad
$ws.Cells.Item($r, 15).Value = "Región del Maule"
$ws.Cells.Item($r, 16).Value = 1353
$ws.Cells.Item($r, 17).Value = $kgUnidades
$ws.Cells.Item($r, 18).Value = $clasif

# New row 666: Morada(o)
$r = 666
$ws.Cells.Item($r, 1).Value  = $mercadoId
$ws.Cells.Item($r, 2).Value  = $mercado
$ws.Cells.Item($r, 3).Value  = $region
$ws.Cells.Item($r, 4).Value  = 44746
$ws.Cells.Item($r, 5).Value  = $codreg
$ws.Cells.Item($r, 6).Value  = $categoriaId
$ws.Cells.Item($r, 7).Value  = $categoria
$ws.Cells.Item($r, 8).Value  = "Morada(o)"
$ws.Cells.Item($r, 9).Value  = $calidad
$ws.Cells.Item($r, 10).Value = 930
$ws.Cells.Item($r, 11).Value = 1300
$ws.Cells.Item($r, 12).Value = 1500
$ws.Cells.Item($r, 13).Value = 1360
$ws.Cells.Item($r, 14).Value = $unidad
$ws.Cells.Item($r, 15).Value = "Región del Maule"
$ws.Cells.Item($r, 16).Value = 1360
$ws.Cells.Item($r, 17).Value = $kgUnidades
$ws.Cells.Item($r, 18).Value = $clasif
